# Update odds values on Sheet1 to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.73

# Row 3
$ws.Range("N3").Value = 8
$ws.Range("Q3").Value = 2.35
$ws.Range("R3").Value = 1.57

# Row 4
$ws.Range("G4").Value = 1.73

# Row 5
$ws.Range("G5").Value = 2.8
$ws.Range("I5").Value = 2.45

# Row 6
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 4.2
$ws.Range("I6").Value = 6.2
$ws.Range("J6").Value = 1.87
$ws.Range("K6").Value = 2.37
$ws.Range("L6").Value = 5.8
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 12.9
$ws.Range("Q6").Value = 1.57
$ws.Range("R6").Value = 2.1
$ws.Range("W6").Value = 6.7
$ws.Range("X6").Value = 6.3
$ws.Range("Z6").Value = 8.25
$ws.Range("AA6").Value = 9.25
$ws.Range("AC6").Value = 13
$ws.Range("AD6").Value = 7.4
$ws.Range("AE6").Value = 14
$ws.Range("AH6").Value = 15
$ws.Range("AI6").Value = 32
$ws.Range("AJ6").Value = 16
$ws.Range("AK6").Value = 100
$ws.Range("AN6").Value = 3.35
$ws.Range("AO6").Value = 6.4
$ws.Range("AQ6").Value = 17.5
$ws.Range("AT6").Value = 3.15
$ws.Range("AU6").Value = 7.7
$ws.Range("AW6").Value = 7.8
$ws.Range("AX6").Value = 35
$ws.Range("AY6").Value = 35

# Row 8
$ws.Range("Q8").Value = 1.85
$ws.Range("R8").Value = 2

# Row 10
$ws.Range("R10").Value = 1.67

# Row 11
$ws.Range("R11").Value = 1.57

# Row 13
$ws.Range("G13").Value = 1.42
